$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values. Each cell is forced to Text format before assignment
# so numeric-looking strings (e.g. '1.005', '0.06615') are preserved as text
# exactly as in the source data, then the style is reset to Normal so no
# extraneous cell style/format is introduced (matching the original workbook
# which left these cells with the default style).
$cellValues = [ordered]@{
    'D2' = '27.895.71'
    'E2' = '  +0.84%  '
    'D3' = '1.756.80'
    'E3' = '  -1.25%  '
    'D4' = '1.005'
    'E4' = '  -0.35%  '
    'D5' = '335.63'
    'E5' = '  -0.87%  '
    'D6' = '1.000'
    'E6' = '  -0.47%  '
    'D7' = '0.3830'
    'E7' = '  -1.70%  '
    'D8' = '0.3387'
    'E8' = '  -1.58%  '
    'D9' = '44.51'
    'E9' = '  -6.79%  '
    'D10' = '1.113'
    'E10' = '  -3.97%  '
    'D11' = '0.07223'
    'E11' = '  -4.32%  '
    'D12' = '1.002'
    'E12' = '  -0.38%  '
    'D13' = '22.46'
    'E13' = '  -1.78%  '
    'D14' = '6.161'
    'E14' = '  -4.58%  '
    'D15' = '7.159'
    'E15' = '  -0.48%  '
    'D16' = '1.757.43'
    'E16' = '  -1.31%  '
    'D17' = '0.00001058'
    'E17' = '  -2.85%  '
    'D18' = '0.06615'
    'E18' = '  -1.43%  '
    'D19' = '79.29'
    'E19' = '  -5.29%  '
    'D20' = '1.001'
    'E20' = '  -0.42%  '
    'B21' = 'Avalanche'
    'C21' = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
    'D21' = '16.66'
    'E21' = '  -6.13%  '
    'B22' = 'Uniswap'
    'C22' = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
    'D22' = '6.226'
    'E22' = '  -4.78%  '
    'D23' = '27.921.22'
    'E23' = '  +0.88%  '
    'D24' = '11.63'
    'E24' = '  -5.74%  '
    'D25' = '2.384'
    'E25' = '  -0.46%  '
    'D26' = '153.04'
    'E26' = '  -1.29%  '
    'D27' = '19.87'
    'E27' = '  -5.83%  '
    'D28' = '2.316'
    'E28' = '  -7.81%  '
    'D29' = '1.959.26'
    'E29' = '  -1.20%  '
    'D30' = '1.287'
    'E30' = '  -13.21%  '
    'D31' = '132.06'
    'E31' = '  -3.72%  '
    'D32' = '4.020'
    'E32' = '  +1.17%  '
    'D33' = '5.833'
    'E33' = '  -7.27%  '
    'D34' = '0.08826'
    'E34' = '  -1.35%  '
    'D35' = '12.24'
    'E35' = '  -6.08%  '
    'D36' = '0.6605'
    'E36' = '  -4.37%  '
    'D37' = '0.06199'
    'E37' = '  -4.20%  '
    'B38' = 'InternetComputer(DFINITY)'
    'C38' = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
    'D38' = '5.155'
    'E38' = '  -6.17%  '
    'B39' = 'VeChain'
    'C39' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    'D39' = '0.02286'
    'E39' = '  -7.38%  '
    'B40' = 'WEMIXTOKEN'
    'C40' = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
    'D40' = '1.520'
    'E40' = '  -3.72%  '
    'D41' = '0.2109'
    'E41' = '  -5.52%  '
    'D42' = '1.210'
    'E42' = '  -3.96%  '
    'D43' = '8.018'
    'E43' = '  -5.76%  '
    'E44' = '  -0.42%  '
    'D45' = '13.72'
    'E45' = '  -4.39%  '
    'D46' = '3.826'
    'E46' = '  -0.97%  '
    'D47' = '0.6046'
    'E47' = '  -5.49%  '
    'D48' = '126.37'
    'E48' = '  -5.43%  '
    'D49' = '2.005'
    'E49' = '  -6.85%  '
    'D50' = '1.132'
    'E50' = '  +5.13%  '
    'D51' = '1.173'
    'E51' = '  +0.42%  '
}

foreach ($cellRef in $cellValues.Keys) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $cellValues[$cellRef]
    $range.Style = "Normal"
}

